## "start course insert into dbf" -- append new vacancy rows (14-18), add a
## "status" column of colour-coded notes (F/G/H depending on row) and
## recolour a few of the rows that already had a verdict (rows 2, 6 -> red
## "declined"/"trash", row 10 -> yellow "passed interview").

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1. New one-off status notes on existing rows
# ---------------------------------------------------------------------
$ws.Range("F2").Value  = "ВЫСЛАЛИ ТЕСТОВОЕ, не та пециальность"
$ws.Range("H6").Value  = "мудаки"
$ws.Range("G10").Value = "Собеседование прошел"

# ---------------------------------------------------------------------
# 2. Four brand new vacancy rows (14-18) in columns A/B/C/D/E, plus a
#    running commentary column G
# ---------------------------------------------------------------------
$ws.Range("A14").Value = "БЭКЭНД-ЭКСПЕРТ"
$ws.Range("G14").Value = "собеседование. Не прошел"

$ws.Range("A15").Value = "VIRON-IT"
$ws.Range("G15").Value = "собеседование, не прошел"

$ws.Range("A16").Value = "Awem Games"
$ws.Range("C16").Value = "https://jobs.tut.by/vacancy/26812835?query=python%20junior"
$ws.Range("E16").Value = "https://vk.com/club16451872"
$ws.Range("G16").Value = "Хм, им питоха не нужен. Но хз, напиши, чо"

$ws.Range("A17").Value = "*instinctools East Europe labs"
$ws.Range("C17").Value = "https://jobs.tut.by/vacancy/26413955?query=Python%20%D1%80%D0%B0%D0%B7%D1%80%D0%B0%D0%B1%D0%BE%D1%82%D1%87%D0%B8%D0%BA"
$ws.Range("D17").Value = "https://www.instinctools.by/?p=597"
$ws.Range("E17").Value = "jobs@instinctools.ru"
$ws.Range("F17").Value = " "

$ws.Range("A18").Value = "Gurtam"
$ws.Range("C18").Value = "https://jobs.tut.by/vacancy/26885408?query=Python%20%D1%80%D0%B0%D0%B7%D1%80%D0%B0%D0%B1%D0%BE%D1%82%D1%87%D0%B8%D0%BA"
$ws.Range("E18").Value = "info@gurtam.com"

# ---------------------------------------------------------------------
# 3. Hyperlinks for the new rows. Adding a hyperlink auto-stamps the
#    "Hyperlink" look onto the cell; re-apply the plain hyperlink cell
#    style (borrowed from an existing hyperlinked cell, C2) straight
#    after so every hyperlinked cell shares the workbook's single
#    hyperlink style instead of growing a fresh one per cell.
# ---------------------------------------------------------------------
$hlStyle = $ws.Range("C2").Style

$ws.Hyperlinks.Add($ws.Range("C16"), "https://jobs.tut.by/vacancy/26812835?query=python%20junior")
$ws.Range("C16").Style = $hlStyle
$ws.Hyperlinks.Add($ws.Range("E16"), "https://vk.com/club16451872")
$ws.Range("E16").Style = $hlStyle
$ws.Hyperlinks.Add($ws.Range("C17"), "https://jobs.tut.by/vacancy/26413955?query=Python%20%D1%80%D0%B0%D0%B7%D1%80%D0%B0%D0%B1%D0%BE%D1%82%D1%87%D0%B8%D0%BA")
$ws.Range("C17").Style = $hlStyle
$ws.Hyperlinks.Add($ws.Range("D17"), "https://www.instinctools.by/?p=597")
$ws.Range("D17").Style = $hlStyle
$ws.Hyperlinks.Add($ws.Range("E17"), "mailto:jobs@instinctools.ru")
$ws.Range("E17").Style = $hlStyle
$ws.Hyperlinks.Add($ws.Range("C18"), "https://jobs.tut.by/vacancy/26885408?query=Python%20%D1%80%D0%B0%D0%B7%D1%80%D0%B0%D0%B1%D0%BE%D1%82%D1%87%D0%B8%D0%BA")
$ws.Range("C18").Style = $hlStyle
$ws.Hyperlinks.Add($ws.Range("E18"), "mailto:info@gurtam.com")
$ws.Range("E18").Style = $hlStyle

# ---------------------------------------------------------------------
# 4. Colour coding -- create the fills in the same order the workbook
#    picked them up (yellow, then red, then blue) so the generated
#    palette lines up with the source file.
# ---------------------------------------------------------------------

# Row 10 -> yellow "passed"
$ws.Range("A10:G10").Interior.Color = 65535

# Row 2 and row 6 -> red "rejected / trash"
$ws.Range("A2:I2").Interior.Color = 255
$ws.Range("A6:H6").Interior.Color = 255

# Scattered single blue marker cells next to a few rows
$ws.Range("F4").Interior.Color  = 15773696
$ws.Range("G5").Interior.Color  = 15773696
$ws.Range("H7").Interior.Color  = 15773696
$ws.Range("I8").Interior.Color  = 15773696
$ws.Range("G9").Interior.Color  = 15773696
$ws.Range("G12").Interior.Color = 15773696
$ws.Range("G13").Interior.Color = 15773696
$ws.Range("G17").Interior.Color = 15773696
$ws.Range("G18").Interior.Color = 15773696

# ---------------------------------------------------------------------
# 5. Selection follows the last row of new data
# ---------------------------------------------------------------------
$ws.Range("G15").Select()
